$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A61").Value = 45965
$ws.Range("A61").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B61").Value = "22,0886"
$ws.Range("C61").Value = "16,1692"
$ws.Range("D61").Value = "15,4027"
$ws.Range("E61").Value = "15,4027"
